$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.55 = 35935.9 pesos`n✅ 35935.9 pesos = 8.52 = 969.71 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update rate cells on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 117
$ws2.Range("O10").Value = 4204.5
$ws2.Range("N12").Value = 4218
$ws2.Range("O12").Value = 113.82
